$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "1605"
$ws.Range("E17").Value = "1606"
$ws.Range("E18").Value = "1607"
$ws.Range("E19").Value = "1608"
$ws.Range("E20").Value = "1609"
$ws.Range("E21").Value = "1610"
$ws.Range("E22").Value = "1611"
$ws.Range("E23").Value = "1612"

$ws.Range("E24").Value = "1807"
$ws.Range("E25").Value = "1808"
